$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '51.095.22'
Set-TextValue 'E2' '  -0.19%  '

Set-TextValue 'D3' '2.961.63'
Set-TextValue 'E3' '  +0.63%  '

Set-TextValue 'D5' '380.58'
Set-TextValue 'E5' '  +1.12%  '

Set-TextValue 'D6' '102.23'
Set-TextValue 'E6' '  -0.38%  '

Set-TextValue 'E7' '  +1.82%  '

Set-TextValue 'E8' '  +0.00%  '

Set-TextValue 'E9' '  +0.69%  '

Set-TextValue 'D10' '36.65'
Set-TextValue 'E10' '  -0.33%  '

Set-TextValue 'E11' '  -0.73%  '

Set-TextValue 'D12' '0.0854'
Set-TextValue 'E12' '  +1.92%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D13' '3.420.84'
Set-TextValue 'E13' '  +0.41%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '18.40'
Set-TextValue 'E14' '  +2.56%  '

Set-TextValue 'D15' '7.75'
Set-TextValue 'E15' '  +5.50%  '

Set-TextValue 'D16' '12.15'
Set-TextValue 'E16' '  +70.87%  '

Set-TextValue 'D17' '2.961.66'
Set-TextValue 'E17' '  +0.63%  '

Set-TextValue 'E18' '  +2.55%  '

Set-TextValue 'D19' '51.179.40'
Set-TextValue 'E19' '  +0.08%  '

Set-TextValue 'D20' '3.11'
Set-TextValue 'E20' '  -1.70%  '

Set-TextValue 'D21' '12.42'
Set-TextValue 'E21' '  -1.36%  '

Set-TextValue 'E22' '  +1.13%  '

Set-TextValue 'D23' '3.34'
Set-TextValue 'E23' '  +15.79%  '

Set-TextValue 'D24' '268.92'
Set-TextValue 'E24' '  +1.92%  '

Set-TextValue 'D25' '69.80'

Set-TextValue 'D26' '7.97'
Set-TextValue 'E26' '  -2.07%  '

Set-TextValue 'E27' '  -0.01%  '

Set-TextValue 'D28' '0.167'
Set-TextValue 'E28' '  -0.87%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D29' '7.11'
Set-TextValue 'E29' '  -9.76%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D30' '25.91'
Set-TextValue 'E30' '  +0.87%  '

Set-TextValue 'E31' '  -3.14%  '

Set-TextValue 'D32' '10.47'
Set-TextValue 'E32' '  +6.19%  '

Set-TextValue 'E33' '  +8.17%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D34' '34.43'
Set-TextValue 'E34' '  +0.59%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D35' '51.27'
Set-TextValue 'E35' '  +0.55%  '

Set-TextValue 'E36' '  -4.14%  '

Set-TextValue 'E37' '  +0.01%  '

Set-TextValue 'D38' '3.29'
Set-TextValue 'E38' '  +10.23%  '

Set-TextValue 'E39' '  +1.86%  '

Set-TextValue 'E40' '  +1.60%  '

Set-TextValue 'E41' '  +3.46%  '

Set-TextValue 'E42' '  -1.76%  '

Set-TextValue 'D43' '124.81'
Set-TextValue 'E43' '  +2.34%  '

Set-TextValue 'E44' '  +10.29%  '

Set-TextValue 'D45' '21.57'
Set-TextValue 'E45' '  +2.50%  '

Set-TextValue 'D46' '0.274'
Set-TextValue 'E46' '  +0.73%  '

Set-TextValue 'E47' '  +3.43%  '

Set-TextValue 'E48' '  -1.62%  '

Set-TextValue 'D49' '2.063.61'
Set-TextValue 'E49' '  +3.22%  '

Set-TextValue 'E50' '  -8.46%  '

Set-TextValue 'E51' '  +6.76%  '
